$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row text. The shared-string table order in the target
# workbook is [delta_h, h1, h2, v_exp], which corresponds to writing the
# cells in this order: D1, B1, C1, E1.
$ws.Range("D1").Value = "delta_h"
$ws.Range("B1").Value = "h1"
$ws.Range("C1").Value = "h2"
$ws.Range("E1").Value = "v_exp"

# New cell H8 (row 8) added below the existing data, with a lightly
# distinguished style (non-default xf referencing the base font).
$ws.Range("H8").Font.ThemeColor = 1
$ws.Range("H8").Select()
